$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G (provincia) metadata rows need re-curating:
# Row2 (dimension qualifier): sdmx-dimension:refArea -> iaest-measure:provincia
$ws.Range("G2").Value = "iaest-measure:provincia"

# Row3 (dim/medida flag): was "dim", should now be "medida" (matching the other measure columns)
$ws.Range("G3").Value = "medida"

# Row4 (datatype): was "URI-Provincia", should now be "xsd:int" (matching the other measure columns)
$ws.Range("G4").Value = "xsd:int"
